$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("linear_increases")

# New notes in column H (rows 38-42) - become shared strings 12-16
$ws.Range("H38").Value = "now I want the total to be 25000"
$ws.Range("H39").Value = "I don't really know what I want it to start at - just not 0"
$ws.Range("H40").Value = "25000 = 1x + 2x +30x….50x"
$ws.Range("H41").Value = "so 25000 = 1275x"
$ws.Range("H42").Value = "x = 20000/1275"

# New constant-manager-budget per-period increment
$ws.Range("H43").Formula = '=5204.1/1275'

# New "manager budget that is constant" note - becomes shared string 17
$ws.Range("L45").Value = "manager budget that is constant"

# F column: time step counters 1..50 for rows 44-93
for ($i = 1; $i -le 50; $i++) {
    $ws.Cells.Item(43 + $i, 6).Value = $i
}

# G column: running manager budget total, starting at 400, incrementing by $H$43
$ws.Range("G44").Value = 400
$ws.Range("G45").Formula = '=G44+$H$43'
$ws.Range("G46").Formula = '=G45+$H$43'
$ws.Range("G47:G93").Formula = '=G46+$H$43'

# Check total sums to (roughly) 25000
$ws.Range("I46").Formula = '=SUM(G44:G93)'

# 25000 divided evenly across the 50 steps
$ws.Range("L46").Formula = '=25000/50'

# Update view to match where editing left off
$ws.Range("H43").Select() | Out-Null
